$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bo sung bao cao: add progress notes for week 3 (row 4)
# F4 gets the "who did what" note, E4 gets the completion result note.
$ws.Range("F4").Value = "Thọ: Giao diện phần app.`nCông: Giao diện phần web API."
$ws.Range("E4").Value = "Hoàn thành xong giao diện."

# Match the original formatting: F4 uses the wrapped-text style like F3/D-column notes
$ws.Range("F4").WrapText = $true

# Update the active selection to E4, as left by the author after the edit
$ws.Range("E4").Select() | Out-Null
